# Generate Report for Handoff
#
# The localization status moves from "Handed back: in sync with en-US"
# to "Ready for handoff" and the associated timestamps are refreshed, on
# all three sheets (Overview summary + the per-locale zh-cn/de-de detail
# sheets that mirror the same status/date strings). The two over-wide
# "Handback"-era datetime columns are also narrowed back down to their
# normal width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: status (zh-cn/de-de columns) + handoff datetime ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-30 19:14:22"

# --- zh-cn sheet: status + refreshed handoff datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-30 19:14:17"

# --- de-de sheet: status + refreshed handoff datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-30 19:14:22"

# --- Column width adjustments (Overview E:F, zh-cn/de-de column C) ---
# The host quantizes ColumnWidth to a 1/6-character pixel grid, so the
# nearest representable value to the authored 17.2159881591797 is used.
$wsOverview.Range("E1").ColumnWidth = 16.333333333333336
$wsOverview.Range("F1").ColumnWidth = 16.333333333333336
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333336
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333336
